$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Neurology" to "Session"
$ws.Name = "Session"

# Remove the last two logged session rows (37 and 38 - the
# 190333/Manual and 191007/Scan entries), shifting remaining cells up so
# the used range shrinks from A1:F38 to A1:F36
$ws.Range("A37:F38").Delete()
